# Generate Report for Handoff
# Update the localization-status report: refresh the "Latest Handoff Datetime"
# / "Latest HO Xliff Generate Date" timestamps and set Priority to "ht" for
# the rows whose handoff xliff files were (re)generated.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-18 20:22:13"
}

# --- zh-cn sheet: Priority (column E) + Latest Handoff Datetime (column H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-18 20:21:59"
}

# --- de-de sheet: Priority (column E) + Latest Handoff Datetime (column H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-18 20:22:13"
}
